# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new value for column F
$exhibitionUpdates = @{
    2  = 24
    7  = 1707
    8  = 31
    11 = 1758
    13 = 105
    21 = 530
    22 = 305
    25 = 258
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new value for column F
$allTypesUpdates = @{
    2  = 24
    7  = 1707
    9  = 31
    12 = 1758
    14 = 105
    22 = 530
    23 = 305
    26 = 258
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
